$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$B_VAL = 33.94444444444444
$C_VAL = 1.95

$rows = @(
    @(0, 0, 0.144),
    @(1, 2, 0.002),
    @(2, 3, 0.005),
    @(3, 4, 0.007),
    @(4, 5, 0.019),
    @(5, 6, 0.033),
    @(6, 7, 0.043),
    @(7, 8, 0.047),
    @(8, 9, 0.044),
    @(9, 10, 0.034),
    @(10, 11, 0.039),
    @(11, 12, 0.034),
    @(12, 13, 0.035),
    @(13, 14, 0.031),
    @(14, 15, 0.04),
    @(15, 16, 0.033),
    @(16, 17, 0.034),
    @(17, 18, 0.031),
    @(18, 19, 0.031),
    @(19, 20, 0.02),
    @(20, 21, 0.024),
    @(21, 22, 0.013),
    @(22, 23, 0.021),
    @(23, 24, 0.023),
    @(24, 25, 0.017),
    @(25, 26, 0.023),
    @(26, 27, 0.021),
    @(27, 28, 0.017),
    @(28, 29, 0.006),
    @(29, 30, 0.011),
    @(30, 31, 0.011),
    @(31, 32, 0.006),
    @(32, 33, 0.016),
    @(33, 34, 0.012),
    @(34, 35, 0.012),
    @(35, 36, 0.009000000000000001),
    @(36, 37, 0.006),
    @(37, 38, 0.006),
    @(38, 39, 0.006),
    @(39, 40, 0.004),
    @(40, 41, 0.005),
    @(41, 42, 0.003),
    @(42, 43, 0.003),
    @(43, 44, 0.004),
    @(44, 45, 0.002),
    @(45, 46, 0.001),
    @(46, 47, 0.001),
    @(47, 48, 0.001),
    @(48, 50, 0.001),
    @(49, 52, 0.002),
    @(50, 53, 0.001),
    @(51, 54, 0.001),
    @(52, 55, 0.001),
    @(53, 56, 0.001),
    @(54, 59, 0.001),
    @(55, 61, 0.001)
)

foreach ($row in $rows) {
    $a = $row[0]
    $d = $row[1]
    $e = $row[2]
    $r = $a + 2
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $B_VAL
    $ws.Cells.Item($r, 3).Value = $C_VAL
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# Copy the style (border/bold/center-top alignment) from A53 to the new A54:A57 cells
$styleSrc = $ws.Cells.Item(53, 1)
$styleSrc.Copy()
for ($r = 54; $r -le 57; $r++) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

Write-Host "Done updating sheet"
